$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 345 previously had "NA" in column C (no page number for the
# catch-all "Rien ne nous concerne aujourd'hui !" entry). The refreshed
# script run leaves that cell blank instead.
$ws.Range("C345").Formula = "'"
$ws.Range("C345").ClearFormats()

# Newly scraped rows appended for 2025-11-27.
$ws.Range("A346").Formula = "'2025-11-27"
$ws.Range("A346").ClearFormats()
$ws.Range("B346").Value = "eaux souterraines"
$ws.Range("C346").Value = 238
$ws.Range("D346").Value = 1

$ws.Range("A347").Formula = "'2025-11-27"
$ws.Range("A347").ClearFormats()
$ws.Range("B347").Value = "eaux souterraines"
$ws.Range("C347").Value = 242
$ws.Range("D347").Value = 2

$ws.Range("A348").Formula = "'2025-11-27"
$ws.Range("A348").ClearFormats()
$ws.Range("B348").Value = "eaux de surface"
$ws.Range("C348").Value = 242
$ws.Range("D348").Value = 1

$ws.Range("A349").Formula = "'2025-11-27"
$ws.Range("A349").ClearFormats()
$ws.Range("B349").Value = "eaux de surface"
$ws.Range("C349").Value = 245
$ws.Range("D349").Value = 1
